$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -2
$ws.Range("F9").Value  = 7
$ws.Range("F10").Value = -1
$ws.Range("F15").Value = 2
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -1
$ws.Range("F22").Value = 1
$ws.Range("F24").Value = -2
$ws.Range("F25").Value = 6
$ws.Range("F30").Value = -1
$ws.Range("F35").Value = -3
$ws.Range("F38").Value = 6
$ws.Range("F39").Value = 8
$ws.Range("F43").Value = 0
$ws.Range("F45").Value = -2
$ws.Range("F46").Value = -8
